$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.189.80"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "2.459.55"
$ws.Range("E3").Value = "  +8.40%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'296.50"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "95.85"
$ws.Range("E6").Value = "  -3.79%  "
$ws.Range("D7").Value = "0.576"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "35.15"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "7.14"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "2.847.55"
$ws.Range("E14").Value = "  +8.77%  "
$ws.Range("D15").Value = "2.461.07"
$ws.Range("E15").Value = "  +8.23%  "
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  +7.66%  "
$ws.Range("D17").Value = "14.17"
$ws.Range("E17").Value = "  +3.83%  "
$ws.Range("D18").Value = "46.267.17"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").Value = "12.71"
$ws.Range("E19").Value = "  +2.48%  "
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  -3.09%  "
$ws.Range("E21").Value = "  +8.23%  "
$ws.Range("D22").Value = "67.57"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").Value = "246.55"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "2.81"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "39.44"
$ws.Range("E27").Value = "  -4.60%  "
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "3.88"
$ws.Range("E30").Value = "  +15.31%  "
$ws.Range("D31").Value = "'21.50"
$ws.Range("E31").Value = "  +7.31%  "
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "5.57"
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("D34").Value = "147.79"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "2.04"
$ws.Range("E35").Value = "  +22.35%  "
$ws.Range("D36").Value = "0.0771"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "15.23"
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("D40").Value = "3.96"
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("D41").Value = "0.0302"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("D42").Value = "3.26"
$ws.Range("E42").Value = "  +6.32%  "
$ws.Range("D43").Value = "1.990.60"
$ws.Range("E43").Value = "  +11.54%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "92.03"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").Value = "16.62"
$ws.Range("E46").Value = "  +34.36%  "
$ws.Range("D47").Value = "1.78"
$ws.Range("E47").Value = "  -6.38%  "
$ws.Range("D48").Value = "'8.60"
$ws.Range("E48").Value = "  +9.37%  "
$ws.Range("D49").Value = "102.08"
$ws.Range("E49").Value = "  +8.06%  "
$ws.Range("D50").Value = "2.712.78"
$ws.Range("E50").Value = "  +8.76%  "
$ws.Range("D51").Value = "0.185"
$ws.Range("E51").Value = "  +0.86%  "
